$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "valor" (D) column with new values ---
$values = @(482.56, 713.21, 619.88, 347.45, 890.33, 250.7, 670.12, 480.94, 790.27, 550.61, 930.18, 360.75, 630.49, 800.02)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $values[$i]
}

# Apply the new number format ("0.00_);[Red]\(0.00\)") to the whole D column data range
$ws.Range("D2:D15").NumberFormat = "0.00_);[Red]\(0.00\)"

# Widen column D to fit the new values
$ws.Columns.Item(4).ColumnWidth = 12.5

# Update the sheet selection to reflect the newly edited column
[void]$ws.Range("D2:D15").Select()

Write-Host "done"
